$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 21.56286567164171
$ws.Range("B2").Value = 1.328
$ws.Range("C2").Value = 45.25599999999999
$ws.Range("A3").Value = 5.004517412935317
$ws.Range("B3").Value = 0.09200000000000004
$ws.Range("C3").Value = 13.99200000000001
$ws.Range("A4").Value = 6.717592039800984
$ws.Range("B4").Value = 0.2959999999999999
$ws.Range("C4").Value = 16.432
$ws.Range("A5").Value = 6.434129353233826
$ws.Range("B5").Value = 0.296
$ws.Range("C5").Value = 17.212
$ws.Range("A6").Value = 13.15367164179102
$ws.Range("B6").Value = 0.9039999999999999
$ws.Range("C6").Value = 29.692
$ws.Range("A7").Value = 19.21826865671642
$ws.Range("B7").Value = 1.016
$ws.Range("C7").Value = 42.76
$ws.Range("A8").Value = 23.25946268656709
$ws.Range("B8").Value = 1.112
$ws.Range("C8").Value = 50.552
$ws.Range("A9").Value = 12.28714427860696
$ws.Range("B9").Value = 0.6800000000000004
$ws.Range("C9").Value = 27.71600000000002
$ws.Range("A10").Value = 26.24220895522376
$ws.Range("B10").Value = 1.308
$ws.Range("C10").Value = 56.83600000000001
$ws.Range("A11").Value = 23.6284975124377
$ws.Range("B11").Value = 1.5
$ws.Range("C11").Value = 50.8
$ws.Range("A12").Value = 18.90330348258704
$ws.Range("B12").Value = 0.8280000000000007
$ws.Range("C12").Value = 42.29600000000001
$ws.Range("A13").Value = 26.11393034825859
$ws.Range("B13").Value = 1.923999999999999
$ws.Range("C13").Value = 57.63200000000002
$ws.Range("A14").Value = 25.49970149253721
$ws.Range("B14").Value = 2.036000000000001
$ws.Range("C14").Value = 54.34799999999998
$ws.Range("A15").Value = 8.065094527363176
$ws.Range("B15").Value = 0.6160000000000003
$ws.Range("C15").Value = 20.324
$ws.Range("A16").Value = 13.75018905472635
$ws.Range("B16").Value = 1.048
$ws.Range("C16").Value = 30.34
$ws.Range("A17").Value = 18.48640796019897
$ws.Range("B17").Value = 1.264
$ws.Range("C17").Value = 40.05600000000001
$ws.Range("A18").Value = 6.786368159203972
$ws.Range("B18").Value = 0.276
$ws.Range("C18").Value = 16.528
$ws.Range("A19").Value = 24.59144278606957
$ws.Range("B19").Value = 1.451999999999999
$ws.Range("C19").Value = 54.07200000000002
$ws.Range("A20").Value = 14.95243781094527
$ws.Range("B20").Value = 0.88
$ws.Range("C20").Value = 32.924
$ws.Range("A21").Value = 10.12644776119403
$ws.Range("B21").Value = 0.5840000000000002
$ws.Range("C21").Value = 24.48400000000002
$ws.Range("A22").Value = 26.82202985074615
$ws.Range("B22").Value = 1.832000000000001
$ws.Range("C22").Value = 56.81600000000001
$ws.Range("A23").Value = 8.030208955223868
$ws.Range("B23").Value = 0.5920000000000001
$ws.Range("C23").Value = 19.936
$ws.Range("A24").Value = 14.11978109452735
$ws.Range("B24").Value = 0.8960000000000002
$ws.Range("C24").Value = 32.06800000000001
$ws.Range("A25").Value = 11.00071641791043
$ws.Range("B25").Value = 0.6040000000000004
$ws.Range("C25").Value = 25.71199999999999
$ws.Range("A26").Value = 9.431343283582084
$ws.Range("B26").Value = 0.82
$ws.Range("C26").Value = 23.488
$ws.Range("A27").Value = 26.40069651741284
$ws.Range("B27").Value = 1.443999999999999
$ws.Range("C27").Value = 57.46
$ws.Range("A28").Value = 20.76827860696507
$ws.Range("B28").Value = 1.12
$ws.Range("C28").Value = 43.55199999999999
$ws.Range("A29").Value = 8.282487562189051
$ws.Range("B29").Value = 0.484
$ws.Range("C29").Value = 20.28399999999999
$ws.Range("A30").Value = 17.61906467661691
$ws.Range("B30").Value = 1.58
$ws.Range("C30").Value = 38.61600000000001
$ws.Range("A31").Value = 23.90803980099493
$ws.Range("B31").Value = 1.284
$ws.Range("C31").Value = 51.20000000000001
$ws.Range("A32").Value = 25.35639800995018
$ws.Range("B32").Value = 1.552
$ws.Range("C32").Value = 55.36400000000001
$ws.Range("A33").Value = 24.95289552238797
$ws.Range("B33").Value = 1.636000000000001
$ws.Range("C33").Value = 54.94399999999997
$ws.Range("A34").Value = 22.43221890547256
$ws.Range("B34").Value = 1.196000000000001
$ws.Range("C34").Value = 48.02400000000002
$ws.Range("A35").Value = 14.25822885572138
$ws.Range("B35").Value = 1.092000000000001
$ws.Range("C35").Value = 32.32400000000003
$ws.Range("A36").Value = 19.72427860696514
$ws.Range("B36").Value = 1.536
$ws.Range("C36").Value = 43.748
$ws.Range("A37").Value = 10.96867661691541
$ws.Range("B37").Value = 0.5040000000000001
$ws.Range("C37").Value = 25.67199999999998
$ws.Range("A38").Value = 26.067422885572
$ws.Range("B38").Value = 1.96
$ws.Range("C38").Value = 53.896
$ws.Range("A39").Value = 19.05832835820895
$ws.Range("B39").Value = 1.02
$ws.Range("C39").Value = 42.47200000000001
$ws.Range("A40").Value = 8.533054726368146
$ws.Range("B40").Value = 0.284
$ws.Range("C40").Value = 21.672
$ws.Range("A41").Value = 17.26189054726367
$ws.Range("B41").Value = 0.9520000000000002
$ws.Range("C41").Value = 39.45600000000002
$ws.Range("A42").Value = 24.03675621890537
$ws.Range("B42").Value = 1.98
$ws.Range("C42").Value = 50.66800000000003
$ws.Range("A43").Value = 17.31016915422885
$ws.Range("B43").Value = 1.427999999999999
$ws.Range("C43").Value = 39.12800000000001
$ws.Range("A44").Value = 20.70819900497496
$ws.Range("B44").Value = 1.196000000000001
$ws.Range("C44").Value = 44.80799999999999
$ws.Range("A45").Value = 24.26718407960188
$ws.Range("B45").Value = 1.532000000000001
$ws.Range("C45").Value = 50.83600000000002
$ws.Range("A46").Value = 23.80790049751232
$ws.Range("B46").Value = 1.48
$ws.Range("C46").Value = 51.21600000000002
$ws.Range("A47").Value = 25.95679601990036
$ws.Range("B47").Value = 1.772
$ws.Range("C47").Value = 58.788
$ws.Range("A48").Value = 19.3397213930348
$ws.Range("B48").Value = 1.244
$ws.Range("C48").Value = 41.732
$ws.Range("A49").Value = 7.98346268656715
$ws.Range("B49").Value = 0.3919999999999999
$ws.Range("C49").Value = 19.92400000000001
$ws.Range("A50").Value = 11.13263681592039
$ws.Range("B50").Value = 0.6599999999999997
$ws.Range("C50").Value = 26.412
$ws.Range("A51").Value = 4.911601990049743
$ws.Range("B51").Value = 0.07600000000000001
$ws.Range("C51").Value = 13.64
$ws.Range("A52").Value = 22.91657711442779
$ws.Range("B52").Value = 1.396
$ws.Range("C52").Value = 49.31199999999998
$ws.Range("A53").Value = 19.09201990049747
$ws.Range("B53").Value = 1.364
$ws.Range("C53").Value = 40.68
$ws.Range("A54").Value = 7.520099502487551
$ws.Range("B54").Value = 0.6720000000000005
$ws.Range("C54").Value = 18.64000000000001
$ws.Range("A55").Value = 10.24326368159203
$ws.Range("B55").Value = 0.8440000000000003
$ws.Range("C55").Value = 23.90000000000001
$ws.Range("A56").Value = 25.09621890547252
$ws.Range("B56").Value = 1.444
$ws.Range("C56").Value = 54.64399999999998
$ws.Range("A57").Value = 25.86907462686556
$ws.Range("B57").Value = 1.740000000000001
$ws.Range("C57").Value = 55.25600000000002
$ws.Range("A58").Value = 10.92232835820894
$ws.Range("B58").Value = 0.6200000000000002
$ws.Range("C58").Value = 25.21599999999999
$ws.Range("A59").Value = 22.37707462686557
$ws.Range("B59").Value = 1.115999999999999
$ws.Range("C59").Value = 49.28000000000001
$ws.Range("A60").Value = 8.597970149253722
$ws.Range("B60").Value = 0.5800000000000003
$ws.Range("C60").Value = 21.87199999999999
$ws.Range("A61").Value = 25.36616915422876
$ws.Range("B61").Value = 1.888
$ws.Range("C61").Value = 54.30399999999999
$ws.Range("A62").Value = 23.48525373134317
$ws.Range("B62").Value = 1.528
$ws.Range("C62").Value = 49.67999999999999
$ws.Range("A63").Value = 25.60107462686553
$ws.Range("B63").Value = 1.504
$ws.Range("C63").Value = 55.64000000000001
$ws.Range("A64").Value = 8.370308457711436
$ws.Range("B64").Value = 0.4920000000000001
$ws.Range("C64").Value = 20.884
$ws.Range("A65").Value = 23.79556218905461
$ws.Range("B65").Value = 1.168
$ws.Range("C65").Value = 51.85199999999999
$ws.Range("A66").Value = 22.8170149253731
$ws.Range("B66").Value = 1.12
$ws.Range("C66").Value = 50.35999999999999
$ws.Range("A67").Value = 22.41717412935317
$ws.Range("B67").Value = 2.131999999999999
$ws.Range("C67").Value = 48.132
$ws.Range("A68").Value = 9.402308457711431
$ws.Range("B68").Value = 0.5320000000000001
$ws.Range("C68").Value = 23.02400000000001
$ws.Range("A69").Value = 26.79474626865663
$ws.Range("B69").Value = 2.004
$ws.Range("C69").Value = 58.31999999999999
$ws.Range("A70").Value = 8.309273631840789
$ws.Range("B70").Value = 0.34
$ws.Range("C70").Value = 20.64799999999999
$ws.Range("A71").Value = 13.07345273631839
$ws.Range("B71").Value = 0.4559999999999998
$ws.Range("C71").Value = 29.74800000000002
$ws.Range("A72").Value = 16.58537313432834
$ws.Range("B72").Value = 1.008
$ws.Range("C72").Value = 35.86399999999998
